$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32; existing rows 32:72 shift down to 33:73.
$ws.Rows.Item(32).EntireRow.Insert()

# Populate the newly inserted row 32 with the new market-day record.
$ws.Range("A32").Value = 3
$ws.Range("B32").Value = "Femacal de La Calera"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = 44789
$ws.Range("E32").Value = 5
$ws.Range("F32").Value = 100112035
$ws.Range("G32").Value = "Bruselas (repollito)"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 40
$ws.Range("K32").Value = 15000
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = 15000
$ws.Range("N32").Value = "`$/malla 15 kilos"
$ws.Range("O32").Value = "Provincia de Quillota"
$ws.Range("P32").Value = 1000
$ws.Range("Q32").Value = 15
$ws.Range("R32").Value = "Hortaliza"
